$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 703, shifting existing rows 703:818 down to 706:821.
$ws.Rows("703:705").Insert()

# Populate the 3 newly-inserted rows (703, 704, 705) with a new weekly entry (date 45127),
# matching the constant structural columns used throughout this sub-table.
foreach ($r in 703..705) {
    $ws.Range("A$r").Value = 11
    $ws.Range("B$r").Value = "Vega Monumental Concepción"
    $ws.Range("C$r").Value = "Bíobío"
    $ws.Range("E$r").Value = 8
    $ws.Range("F$r").Value = "Fruta"
    $ws.Range("G$r").Value = 100108
    $ws.Range("H$r").Value = "Tropicales y subtropicales"
    $ws.Range("I$r").Value = 100108006
    $ws.Range("J$r").Value = "Plátano"
    $ws.Range("K$r").Value = "Sin especificar"
    $ws.Range("Q$r").Value = "`$/caja 20 kilos"
    $ws.Range("R$r").Value = "Ecuador"
    $ws.Range("T$r").Value = 20
}

# Row 703: Maduro
$ws.Range("D703").Value = 45127
$ws.Range("L703").Value = "Maduro"
$ws.Range("M703").Value = 100
$ws.Range("N703").Value = 11000
$ws.Range("O703").Value = 11000
$ws.Range("P703").Value = 11000
$ws.Range("S703").Value = 550

# Row 704: Pintón
$ws.Range("D704").Value = 45127
$ws.Range("L704").Value = "Pintón"
$ws.Range("M704").Value = 400
$ws.Range("N704").Value = 12000
$ws.Range("O704").Value = 12000
$ws.Range("P704").Value = 12000
$ws.Range("S704").Value = 600

# Row 705: Primera Pintón
$ws.Range("D705").Value = 45127
$ws.Range("L705").Value = "Primera Pintón"
$ws.Range("M705").Value = 400
$ws.Range("N705").Value = 15000
$ws.Range("O705").Value = 15000
$ws.Range("P705").Value = 15000
$ws.Range("S705").Value = 750
